$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.771.05'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '1.630.53'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.71%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.53'
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.500'
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("E7").Value = '  -0.71%  '
$ws.Range("E8").Value = '  -0.73%  '
$ws.Range("E9").Value = '  -0.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.55'
$ws.Range("E10").Value = '  +1.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0792'
$ws.Range("E11").Value = '  +1.12%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.674.04'
$ws.Range("E12").Value = '  +2.68%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.26'
$ws.Range("E13").Value = '  +0.64%  '
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '1.856.68'
$ws.Range("E14").Value = '  -0.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.554'
$ws.Range("E15").Value = '  +0.78%  '
$ws.Range("E16").Value = '  -0.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.91'
$ws.Range("E17").Value = '  -0.34%  '
$ws.Range("D18").Value = '25.756.15'
$ws.Range("E18").Value = '  +0.01%  '
$ws.Range("E19").Value = '  -0.64%  '
$ws.Range("E20").Value = '  +0.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '191.34'
$ws.Range("E21").Value = '  -0.67%  '
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.27'
$ws.Range("E23").Value = '  +1.43%  '
$ws.Range("E24").Value = '  +2.09%  '
$ws.Range("E25").Value = '  -0.72%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.49'
$ws.Range("E26").Value = '  +1.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.124'
$ws.Range("E27").Value = '  +4.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.84'
$ws.Range("E28").Value = '  +0.60%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  -0.23%  '
$ws.Range("E31").Value = '  +2.08%  '
$ws.Range("E32").Value = '  +0.23%  '
$ws.Range("E33").Value = '  -0.18%  '
$ws.Range("E34").Value = '  +0.55%  '
$ws.Range("E35").Value = '  -0.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.907'
$ws.Range("E36").Value = '  +1.52%  '
$ws.Range("D37").Value = '1.137.98'
$ws.Range("E37").Value = '  +2.96%  '
$ws.Range("E38").Value = '  -2.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.542'
$ws.Range("E39").Value = '  -0.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0155'
$ws.Range("E40").Value = '  -0.37%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.998'
$ws.Range("E41").Value = '  -0.68%  '
$ws.Range("E42").Value = '  -0.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.77'
$ws.Range("E43").Value = '  +1.30%  '
$ws.Range("E44").Value = '  -0.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.803'
$ws.Range("E45").Value = '  +0.35%  '
$ws.Range("D46").Value = '1.765.50'
$ws.Range("E46").Value = '  +0.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.18'
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0510'
$ws.Range("E48").Value = '  +1.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.44'
$ws.Range("E49").Value = '  +6.14%  '
$ws.Range("E50").Value = '  -0.21%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.34'
$ws.Range("E51").Value = '  -0.51%  '
